$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($worksheet, $row, $values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $worksheet.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Row 2 (existing row): Target cluster changes from sCs -> M2, metrics recomputed
Set-RowValues $ws 2 @(
    "sCs", "Bdnf", "Ngfr", "M2",
    3, 1, 4.674406333333334, 14.023219, 1, 1, 2, 0.6666666666666666,
    0.2652796666666666, 0.795839, 0.03859763870461846, 0.03859763870461847,
    1.240024953971222, 11.160224585741, 0.03859763870461846, 0.03859763870461847
)

# Row 3 (new row): Target cluster = Neutro
Set-RowValues $ws 3 @(
    "sCs", "Bdnf", "Ngfr", "Neutro",
    3, 1, 4.674406333333334, 14.023219, 1, 1, 1, 0.3333333333333333,
    0.034619, 0.103857, 0.005036992360195417, 0.005036992360195417,
    0.1618232728536667, 1.456409455683, 0.005036992360195417, 0.005036992360195417
)

# Row 4 (new row): Target cluster = sCs
Set-RowValues $ws 4 @(
    "sCs", "Bdnf", "Ngfr", "sCs",
    3, 1, 4.674406333333334, 14.023219, 1, 1, 3, 1,
    6.573052000000001, 19.719156, 0.9563653689351861, 0.9563653689351862,
    30.72511589812934, 276.526043083164, 0.9563653689351861, 0.9563653689351862
)
